# Auto-generated script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.917.52"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.810.66"
$ws.Range("E3").Value = "  +2.09%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.32%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.13"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +1.17%  "
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.03%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4285"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -2.29%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3691"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +1.21%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07246"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +0.72%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8615"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +2.85%  "
$ws.Range("D11").Value = "2.033.54"
$ws.Range("E11").Value = "  +15.86%  "
$ws.Range("E12").Value = "  +4.58%  "
$ws.Range("E13").Value = "  +4.41%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.384"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +2.50%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06896"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +1.48%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "80.63"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("E17").Value = "  -0.33%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008911"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("E19").Value = "  +0.01%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.16"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "26.936.10"
$ws.Range("E21").Value = "  +1.17%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.189"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +3.40%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.94"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").Value = "2.251.77"
$ws.Range("E24").Value = "  +14.15%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.69"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +0.29%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.884"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -1.10%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.31"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  +3.21%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.896"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +15.44%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.02"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +0.58%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08930"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -0.78%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7420"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +3.30%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +6.46%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.416"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +2.19%  "
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.801"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  +0.06%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("E38").Value = "  +2.16%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01918"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +1.64%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5074"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +3.06%  "
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.732"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +6.28%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1643"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +2.22%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.423"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +4.97%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.238"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +4.01%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.92"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +2.02%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +3.27%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.648"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +4.56%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06278"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4549"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +1.38%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.796"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +5.40%  "
